$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 26 needs the same style (s="1", bold/bordered/centered) as the
# other label cells in column A -- copy format+value from A2, then
# overwrite the value below.
$ws.Range("A2").Copy($ws.Range("A26"))

# Update column A labels (rows 2-26), reflecting the re-ordering / new model row
$ws.Cells.Item(2, 1).Value = "model_2_5_0"
$ws.Cells.Item(3, 1).Value = "model_2_5_22"
$ws.Cells.Item(4, 1).Value = "model_2_5_21"
$ws.Cells.Item(5, 1).Value = "model_2_5_20"
$ws.Cells.Item(6, 1).Value = "model_2_5_19"
$ws.Cells.Item(7, 1).Value = "model_2_5_18"
$ws.Cells.Item(8, 1).Value = "model_2_5_17"
$ws.Cells.Item(9, 1).Value = "model_2_5_16"
$ws.Cells.Item(10, 1).Value = "model_2_5_15"
$ws.Cells.Item(11, 1).Value = "model_2_5_14"
$ws.Cells.Item(12, 1).Value = "model_2_5_13"
$ws.Cells.Item(13, 1).Value = "model_2_5_23"
$ws.Cells.Item(14, 1).Value = "model_2_5_12"
$ws.Cells.Item(15, 1).Value = "model_2_5_10"
$ws.Cells.Item(16, 1).Value = "model_2_5_9"
$ws.Cells.Item(17, 1).Value = "model_2_5_8"
$ws.Cells.Item(18, 1).Value = "model_2_5_7"
$ws.Cells.Item(19, 1).Value = "model_2_5_6"
$ws.Cells.Item(20, 1).Value = "model_2_5_5"
$ws.Cells.Item(21, 1).Value = "model_2_5_4"
$ws.Cells.Item(22, 1).Value = "model_2_5_3"
$ws.Cells.Item(23, 1).Value = "model_2_5_2"
$ws.Cells.Item(24, 1).Value = "model_2_5_1"
$ws.Cells.Item(25, 1).Value = "model_2_5_11"
$ws.Cells.Item(26, 1).Value = "model_2_5_24"

# Update columns B..I (metrics) for rows 2-26 -- all rows now share model_2_5_24 metrics
$metricValues = @(0.85251841622735, 0.7579013496919393, 0.9158305375510991, 0.8422238908847012, 0.1632186621427536, 0.2474200129508972, 0.1003445982933044, 0.1782080680131912)
for ($r = 2; $r -le 26; $r++) {
    for ($i = 0; $i -lt 8; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $metricValues[$i]
    }
}

